$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 6
$ws.Range("I2").Value = 1.62
$ws.Range("J2").Value = 6.5
$ws.Range("L2").Value = 2.3
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 11
$ws.Range("Y2").Value = 19
$ws.Range("AB2").Value = 51
$ws.Range("AI2").Value = 6.5
$ws.Range("AL2").Value = 17
$ws.Range("AO2").Value = 34
$ws.Range("AW2").Value = 3.4
$ws.Range("AY2").Value = 23

# Row 4 updates
$ws.Range("K4").Value = 2.25
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.2
$ws.Range("P4").Value = 4.33
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 2.1
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 12
$ws.Range("AG4").Value = 151
$ws.Range("AJ4").Value = 10
$ws.Range("AK4").Value = 26
$ws.Range("AO4").Value = 15
$ws.Range("AU4").Value = 7.5
$ws.Range("BD4").Value = 126
